$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph "ligacao(...)" -- drop the spellStart/spellEnd proofing
# marks that wrapped "meio_transporte" (the word is no longer flagged).
# ---------------------------------------------------------------------------
$p24 = $d.Paragraphs.Item(24).Range
if ($p24.Text -notlike "liga*o(local_a, local_b, meio_transporte, dist*") {
    throw "Paragraph 24 did not match the expected 'ligacao(...)' text: $($p24.Text)"
}

$p24xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00A56285" w:rsidRDefault="00A56285" w:rsidP="006A04B8"><w:pPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>liga&#231;&#227;o</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidR="00871440"><w:rPr><w:sz w:val="28"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr><w:t>local_</w:t></w:r><w:r w:rsidR="009D0BC6" w:rsidRPr="00EB0ACC"><w:rPr><w:sz w:val="28"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">a, local_b, </w:t></w:r><w:r w:rsidR="00871440"><w:rPr><w:sz w:val="28"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr><w:t>meio_</w:t></w:r><w:r w:rsidR="00EB0ACC" w:rsidRPr="00EB0ACC"><w:rPr><w:sz w:val="28"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr><w:t>transporte</w:t></w:r><w:r w:rsidR="00EB0ACC"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>,</w:t></w:r><w:r w:rsidR="001C664D"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> dist&#226;ncia,</w:t></w:r><w:r w:rsidR="00EB0ACC"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00871440"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>informa&#231;&#227;o_extra</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00871440"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00871440"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>ponto_cardeal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00290150"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00290150"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>ordem_a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00290150"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00290150"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>ordem_b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00871440"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@

$p24.InsertXML($p24xml)

# ---------------------------------------------------------------------------
# Change 2: the "destino(...)" definition paragraph plus the three rule
# paragraphs that followed it are reworked:
#   * destino(...) gets a new, reordered attribute list with "destino" (not
#     "origem, sentido, destino") as the underlined primary-key part, and the
#     trailing "_GoBack" bookmark now sits at the end of *this* paragraph.
#   * the "local_a, local_b, meio_transporte: FK(ligacao)" paragraph is kept
#     but re-tagged with gramStart/gramEnd proofing marks.
#   * RI-13's wording changes completely.
#   * RI-14 is removed outright.
# ---------------------------------------------------------------------------
$first = $d.Paragraphs.Item(31).Range
$last = $d.Paragraphs.Item(34).Range
if ($first.Text -notlike "destino(origem, sentido, destino*") {
    throw "Paragraph 31 did not match the expected 'destino(...)' text: $($first.Text)"
}
if ($last.Text -notlike "RI-14*") {
    throw "Paragraph 34 did not match the expected 'RI-14' text: $($last.Text)"
}

$block = $d.Range($first.Start, $last.End)

$blockxml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000E3AF1" w:rsidRDefault="00582392" w:rsidP="000E3AF1"><w:pPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>destino</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidR="00910A78" w:rsidRPr="00910A78"><w:rPr><w:sz w:val="28"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">local_a, local_b, meio_transporte, origem, </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr><w:t>destino</w:t></w:r><w:r w:rsidR="00910A78"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="009F5B0B" w:rsidRDefault="009F5B0B" w:rsidP="000E3AF1"><w:pPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>local_a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">, local_b, meio_transporte: </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>FK(liga&#231;&#227;o</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p w:rsidR="00F30EB6" w:rsidRDefault="00DB7FC8" w:rsidP="000E3AF1"><w:pPr><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>RI-13</w:t></w:r><w:r w:rsidR="00F30EB6"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">: A </w:t></w:r><w:r w:rsidR="00F30EB6"><w:rPr><w:i/><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t>origem</w:t></w:r><w:r w:rsidR="00F30EB6"><w:rPr><w:sz w:val="28"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> s&#243; pode ter um de 2 valores diferentes</w:t></w:r></w:p>
'@

$block.InsertXML($blockxml)
